$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101:105 down to 102:106
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new weekly record
$ws.Range("A101").Value = 5
$ws.Range("B101").Value = "Macroferia Regional de Talca"
$ws.Range("C101").Value = "Maule"
$ws.Range("D101").Value = 44610
$ws.Range("E101").Value = 7
$ws.Range("F101").Value = 100112001
$ws.Range("G101").Value = "Berenjena"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 200
$ws.Range("K101").Value = 7000
$ws.Range("L101").Value = 7000
$ws.Range("M101").Value = 7000
$ws.Range("N101").Value = "`$/caja 50 unidades"
$ws.Range("O101").Value = "Región del Maule"
$ws.Range("P101").Value = 140
$ws.Range("Q101").Value = 50
$ws.Range("R101").Value = "Hortaliza"
